# Generate Report for Handoff
# Adds a new tracked file (da24845a-cfbd-4a14-ba6c-ef09d047a7c8.md) as row 9
# on the Overview / zh-cn / de-de sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$uuid        = "da24845a-cfbd-4a14-ba6c-ef09d047a7c8"
$mdName      = "$uuid.md"
$zhXlfName   = "$uuid.dc753fc5d951b7e18cb7a551ae77846fb46268c1.zh-cn.xlf"
$deXlfName   = "$uuid.dc753fc5d951b7e18cb7a551ae77846fb46268c1.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1): File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"
$wsOverview.Range("D9").Value = "2016-33-17 16:33:01"

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2): Source File Name | File Extension | Status |
#   Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#   Latest Handback File | Latest Handback DateTime | Handoff Reason |
#   Dependency From | Error Detail
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    ".md"
)
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName",
    "",
    "",
    $zhXlfName
)
$wsZhCn.Range("E9").Value = "2016-03-17 16:32:57"
$wsZhCn.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I9").Value = "Include"

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3): same columns as zh-cn
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    $mdName
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0000000000000000000000000000000000000000/e2e/$mdName",
    "",
    "",
    ".md"
)
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName",
    "",
    "",
    $deXlfName
)
$wsDeDe.Range("E9").Value = "2016-03-17 16:33:01"
$wsDeDe.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I9").Value = "Include"

Write-Host "Added handoff row for $mdName to Overview, zh-cn, de-de sheets."
